# Updated symbol list on Thu Feb  2 22:32:07 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) quotes for the
# coinranking.com crypto price sheet. Values are stored as plain text in
# the source workbook, so each assignment is apostrophe-prefixed to force
# Excel to keep them as text (matching the original inlineStr cells)
# instead of silently re-typing them as Number/Percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.01"
$ws.Range("E2").Value = "'3.39%"
$ws.Range("D3").Value = "'40.18"
$ws.Range("E3").Value = "'5.55%"
$ws.Range("D4").Value = "'5.806"
$ws.Range("E4").Value = "'12.37%"
$ws.Range("D5").Value = "'0.08027"
$ws.Range("E5").Value = "'0.50%"
$ws.Range("D6").Value = "'4.571"
$ws.Range("E6").Value = "'2.19%"
$ws.Range("D7").Value = "'8.709"
$ws.Range("E7").Value = "'2.14%"
$ws.Range("D8").Value = "'1.949"
$ws.Range("E8").Value = "'0.70%"
$ws.Range("E9").Value = "'-0.43%"
$ws.Range("D10").Value = "'0.9450"
$ws.Range("E10").Value = "'0.24%"
$ws.Range("D11").Value = "'0.1250"
$ws.Range("E11").Value = "'-2.36%"
$ws.Range("D12").Value = "'0.1956"
$ws.Range("E12").Value = "'0.78%"
$ws.Range("D13").Value = "'8.923"
$ws.Range("E13").Value = "'35.40%"
$ws.Range("D14").Value = "'0.09228"
$ws.Range("E14").Value = "'1.75%"
$ws.Range("D15").Value = "'0.03584"
$ws.Range("E15").Value = "'4.94%"
$ws.Range("D16").Value = "'0.09637"
$ws.Range("E16").Value = "'1.24%"
$ws.Range("D17").Value = "'0.001300"
$ws.Range("E17").Value = "'-4.37%"
$ws.Range("D18").Value = "'0.006157"
$ws.Range("E18").Value = "'0.80%"
$ws.Range("D19").Value = "'3.369"
$ws.Range("E19").Value = "'-1.69%"
$ws.Range("D20").Value = "'0.3526"
$ws.Range("E20").Value = "'0.33%"
$ws.Range("D21").Value = "'0.1405"
$ws.Range("E21").Value = "'7.64%"
$ws.Range("D22").Value = "'0.2418"
$ws.Range("E22").Value = "'0.11%"
$ws.Range("D23").Value = "'0.04409"
$ws.Range("E23").Value = "'0.82%"
$ws.Range("D24").Value = "'0.001262"
$ws.Range("E24").Value = "'2.93%"
$ws.Range("D25").Value = "'0.004326"
$ws.Range("E25").Value = "'1.46%"
$ws.Range("D26").Value = "'0.0001146"
$ws.Range("E27").Value = "'0.53%"
$ws.Range("D39").Value = "'0.02423"
$ws.Range("E39").Value = "'-0.79%"
$ws.Range("D40").Value = "'0.05262"
$ws.Range("E40").Value = "'2.15%"
$ws.Range("D41").Value = "'0.007463"
$ws.Range("E41").Value = "'-2.30%"
$ws.Range("D42").Value = "'0.1419"
$ws.Range("E42").Value = "'1.30%"
$ws.Range("D43").Value = "'0.008557"
$ws.Range("E43").Value = "'0.30%"
$ws.Range("D44").Value = "'0.002110"
$ws.Range("E44").Value = "'0.21%"
$ws.Range("D45").Value = "'0.01094"
$ws.Range("E45").Value = "'25.15%"
$ws.Range("D46").Value = "'0.00006904"
$ws.Range("E46").Value = "'6.88%"
$ws.Range("E47").Value = "'0.64%"
$ws.Range("D48").Value = "'0.003160"
$ws.Range("E48").Value = "'10.48%"
$ws.Range("D49").Value = "'0.001426"
$ws.Range("E49").Value = "'-15.20%"
$ws.Range("E50").Value = "'0.64%"
$ws.Range("E51").Value = "'0.64%"
